$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1. valid_product_names: remove the 16 "life product" rows that used
#    to occupy Nigeria's column (D2:D17), shifting the remaining
#    Nigeria product names up so they start again at D2.
# -----------------------------------------------------------------
$wsValid = $wb.Worksheets.Item("valid_product_names")

# Capture the Nigeria values that should survive (previously D18:D54)
$nigeriaVals = @()
for ($r = 18; $r -le 54; $r++) {
    $nigeriaVals += $wsValid.Cells.Item($r, 4).Value2
}

# Clear out the whole existing Nigeria column data range
$wsValid.Range("D2:D54").ClearContents()

# Write the surviving values back starting at D2 (i.e. shift up by 16)
$row = 2
foreach ($val in $nigeriaVals) {
    $wsValid.Cells.Item($row, 4).Value2 = $val
    $row++
}

# -----------------------------------------------------------------
# 2. Update the defined names so each country range matches the new
#    extents of data in valid_product_names.
# -----------------------------------------------------------------
$wb.Names.Item("Kenya").RefersTo   = "=valid_product_names!`$A`$2:`$A`$62"
$wb.Names.Item("Nigeria").RefersTo = "=valid_product_names!`$D`$2:`$D`$38"
$wb.Names.Item("Uganda").RefersTo  = "=valid_product_names!`$C`$2:`$C`$51"
$wb.Names.Item("Zimbabwe").RefersTo = "=valid_product_names!`$B`$2:`$B`$62"

# -----------------------------------------------------------------
# 3. customer_prod: refresh the sample/demo row so it points at a
#    Nigeria product instead of the now-pruned Ghana example.
# -----------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("customer_prod")
$wsProd.Range("C2").Value2 = "Nigeria"
$wsProd.Range("A2").Value2 = "Specified Cars"
$wsProd.Range("A3").Value2 = "Householders"

# -----------------------------------------------------------------
# 4. customer_acc: add two more example account numbers to the
#    sample list.
# -----------------------------------------------------------------
$wsAcc = $wb.Worksheets.Item("customer_acc")
$wsAcc.Range("A12").Value2 = "3013483433"
$wsAcc.Range("A13").Value2 = "3013460492"
